$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (D1, style index 1) across the new header range E1:AS1
$ws.Range("D1").Copy()
$ws.Range("E1:AS1").PasteSpecial(-4122)

# Set header text for each column (D1 is retitled; E1:AS1 are new columns)
$ws.Range("D1").Value = "Study ID"
$ws.Range("E1").Value = "Country/Region"
$ws.Range("F1").Value = "Study Design"
$ws.Range("G1").Value = "Database/Setting"
$ws.Range("H1").Value = "Age (mean ± SD)"
$ws.Range("I1").Value = "Sample Size (Total)"
$ws.Range("J1").Value = "GLP-1 RA Cohort Size"
$ws.Range("K1").Value = "Control cohort Size"
$ws.Range("L1").Value = "Sex (% male/female)"
$ws.Range("M1").Value = "Diabetes Status (%)"
$ws.Range("N1").Value = "Obesity/BMI (mean or %)"
$ws.Range("O1").Value = "Smoking History (%)"
$ws.Range("P1").Value = "Comorbidities"
$ws.Range("Q1").Value = "Indication for Spine Surgery"
$ws.Range("R1").Value = "Surgical Procedure"
$ws.Range("S1").Value = "Levels of Surgery"
$ws.Range("T1").Value = "Follow-up Duration"
$ws.Range("U1").Value = "Bone Mineral Density (BMD)"
$ws.Range("V1").Value = "Glycemic Control (HbA1c)"
$ws.Range("W1").Value = "Glucagon-like peptide-1 receptor agonist Details"
$ws.Range("X1").Value = "Preoperative GLP-1 RA exposure duration"
$ws.Range("Y1").Value = "Control Group Details"
$ws.Range("Z1").Value = "Primary Outcome(s)"
$ws.Range("AA1").Value = "Secondary Outcome(s)"
$ws.Range("AB1").Value = "Intraoperative Complications"
$ws.Range("AC1").Value = "Postoperative Complications"
$ws.Range("AD1").Value = "Surgical Site Infection (SSI)"
$ws.Range("AE1").Value = "Reoperation Rates"
$ws.Range("AF1").Value = "Pseudoarthrosis/Fusion Failure"
$ws.Range("AG1").Value = "Readmission"
$ws.Range("AH1").Value = "Mortality"
$ws.Range("AI1").Value = "Other Medical Complications"
$ws.Range("AJ1").Value = "Blood Loss"
$ws.Range("AK1").Value = "Operative Time"
$ws.Range("AL1").Value = "Hospital Length of Stay (LOS)"
$ws.Range("AM1").Value = "Patient-Reported Outcomes (PROMs)"
$ws.Range("AN1").Value = "Cost/Economic Analysis"
$ws.Range("AO1").Value = "Key Findings/Conclusion"
$ws.Range("AP1").Value = "Subgroup Analysis"
$ws.Range("AQ1").Value = "Heterogeneity"
$ws.Range("AR1").Value = "Risk of Bias/Quality Assessment"
$ws.Range("AS1").Value = "Source File"

Write-Output "DONE"
